$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.305.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.68%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.854.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.25%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4557'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.68%  '

$ws.Range("E8").Value = '  -2.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.21'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07912'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.012'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.92%  '

$ws.Range("E12").Value = '  -3.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.866.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.911'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.163'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.87%  '

$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06633'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '86.04'
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '  -3.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.70%  '

$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.501'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.314.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.285'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.076.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.062'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.463'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9454'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09353'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.444'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.593'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.256'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06041'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.71%  '

$ws.Range("E38").Value = '  -3.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.217'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.057'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.28%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5927'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.10%  '

$ws.Range("E43").Value = '  -0.71%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.46%  '

$ws.Range("E45").Value = '  -1.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5620'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.380'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.914'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06740'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '108.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.98%  '
